# Update the cryptos price/volume table with freshly scraped values.
# D column holds "Price" (stored as text, since source values use a
# non-numeric dotted-thousands style like "25.844.56"); E column holds
# "Volume(1h)" percentages (always text, e.g. "  +0.22%  ").
#
# For D-column values that look like plain numbers (e.g. "1.000",
# "0.9996") we prefix the assignment with a leading apostrophe so Excel
# keeps them as text instead of silently converting them to numeric
# values, then reset the cell Style back to "Normal" so no stray
# quote-prefix formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.844.56'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '1.738.35'
$ws.Range('E3').Value = '  -0.85%  '
$r = $ws.Range('D4')
$r.Value = "'0.9996"
$r.Style = 'Normal'
$ws.Range('E4').Value = '  -0.11%  '
$r = $ws.Range('D5')
$r.Value = "'227.48"
$r.Style = 'Normal'
$ws.Range('E5').Value = '  -4.11%  '
$r = $ws.Range('D6')
$r.Value = "'1.000"
$r.Style = 'Normal'
$ws.Range('E6').Value = '  -0.03%  '
$r = $ws.Range('D7')
$r.Value = "'0.5147"
$r.Style = 'Normal'
$ws.Range('E7').Value = '  +1.36%  '
$r = $ws.Range('D8')
$r.Value = "'0.2683"
$r.Style = 'Normal'
$ws.Range('E8').Value = '  +1.33%  '
$r = $ws.Range('D9')
$r.Value = "'39.33"
$r.Style = 'Normal'
$ws.Range('E9').Value = '  -5.29%  '
$r = $ws.Range('D10')
$r.Value = "'0.06083"
$r.Style = 'Normal'
$ws.Range('E10').Value = '  -1.39%  '
$ws.Range('D11').Value = '1.734.84'
$ws.Range('E11').Value = '  -1.23%  '
$r = $ws.Range('D12')
$r.Value = "'0.06997"
$r.Style = 'Normal'
$ws.Range('E12').Value = '  +1.13%  '
$r = $ws.Range('D13')
$r.Value = "'15.18"
$r.Style = 'Normal'
$ws.Range('E13').Value = '  -3.40%  '
$r = $ws.Range('D14')
$r.Value = "'0.6272"
$r.Style = 'Normal'
$ws.Range('E14').Value = '  +3.87%  '
$r = $ws.Range('D15')
$r.Value = "'4.485"
$r.Style = 'Normal'
$ws.Range('E15').Value = '  -0.46%  '
$r = $ws.Range('D16')
$r.Value = "'76.33"
$r.Style = 'Normal'
$ws.Range('E16').Value = '  -1.18%  '
$r = $ws.Range('D17')
$r.Value = "'0.9986"
$r.Style = 'Normal'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('E18').Value = '  -0.03%  '
$ws.Range('D19').Value = '25.850.88'
$ws.Range('E19').Value = '  +0.24%  '
$r = $ws.Range('D20')
$r.Value = "'11.42"
$r.Style = 'Normal'
$ws.Range('E20').Value = '  -2.53%  '
$r = $ws.Range('D21')
$r.Value = "'0.000006531"
$r.Style = 'Normal'
$ws.Range('E21').Value = '  -4.77%  '
$ws.Range('D22').Value = '1.957.76'
$ws.Range('E22').Value = '  -0.97%  '
$r = $ws.Range('D23')
$r.Value = "'4.021"
$r.Style = 'Normal'
$ws.Range('E23').Value = '  -1.65%  '
$r = $ws.Range('D24')
$r.Value = "'8.356"
$r.Style = 'Normal'
$ws.Range('E24').Value = '  +1.14%  '
$r = $ws.Range('D25')
$r.Value = "'5.084"
$r.Style = 'Normal'
$ws.Range('E25').Value = '  -2.50%  '
$r = $ws.Range('D26')
$r.Value = "'136.64"
$r.Style = 'Normal'
$ws.Range('E26').Value = '  -0.73%  '
$r = $ws.Range('D27')
$r.Value = "'1.506"
$r.Style = 'Normal'
$ws.Range('E27').Value = '  +2.42%  '
$ws.Range('E28').Value = '  -0.88%  '
$r = $ws.Range('D29')
$r.Value = "'14.95"
$r.Style = 'Normal'
$ws.Range('E29').Value = '  -0.56%  '
$r = $ws.Range('D30')
$r.Value = "'102.78"
$r.Style = 'Normal'
$ws.Range('E30').Value = '  -0.10%  '
$r = $ws.Range('D31')
$r.Value = "'0.08290"
$r.Style = 'Normal'
$ws.Range('E31').Value = '  +0.90%  '
$r = $ws.Range('D32')
$r.Value = "'3.611"
$r.Style = 'Normal'
$ws.Range('E32').Value = '  -2.19%  '
$r = $ws.Range('D33')
$r.Value = "'3.344"
$r.Style = 'Normal'
$ws.Range('E33').Value = '  -3.57%  '
$r = $ws.Range('D34')
$r.Value = "'0.04411"
$r.Style = 'Normal'
$ws.Range('E34').Value = '  -2.40%  '
$r = $ws.Range('D35')
$r.Value = "'2.612"
$r.Style = 'Normal'
$ws.Range('E35').Value = '  -1.76%  '
$r = $ws.Range('D36')
$r.Value = "'0.9740"
$r.Style = 'Normal'
$ws.Range('E36').Value = '  -2.77%  '
$r = $ws.Range('D37')
$r.Value = "'0.5947"
$r.Style = 'Normal'
$ws.Range('E37').Value = '  -2.25%  '
$r = $ws.Range('D38')
$r.Value = "'2.676"
$r.Style = 'Normal'
$ws.Range('E38').Value = '  -0.80%  '
$r = $ws.Range('D39')
$r.Value = "'0.01565"
$r.Style = 'Normal'
$ws.Range('E39').Value = '  +0.49%  '
$r = $ws.Range('D40')
$r.Value = "'1.914"
$r.Style = 'Normal'
$ws.Range('E40').Value = '  -1.94%  '
$r = $ws.Range('D41')
$r.Value = "'0.9987"
$r.Style = 'Normal'
$ws.Range('E41').Value = '  -0.17%  '
$r = $ws.Range('D42')
$r.Value = "'101.76"
$r.Style = 'Normal'
$ws.Range('E42').Value = '  -1.84%  '
$r = $ws.Range('D43')
$r.Value = "'0.3791"
$r.Style = 'Normal'
$ws.Range('E43').Value = '  -1.12%  '
$r = $ws.Range('D44')
$r.Value = "'0.7284"
$r.Style = 'Normal'
$ws.Range('E44').Value = '  -1.68%  '
$r = $ws.Range('D45')
$r.Value = "'4.814"
$r.Style = 'Normal'
$ws.Range('E45').Value = '  -2.23%  '
$ws.Range('E46').Value = '  +0.04%  '
$r = $ws.Range('D47')
$r.Value = "'6.218"
$r.Style = 'Normal'
$ws.Range('E47').Value = '  +3.46%  '
$r = $ws.Range('D48')
$r.Value = "'0.1100"
$r.Style = 'Normal'
$ws.Range('E48').Value = '  -0.32%  '
$r = $ws.Range('D49')
$r.Value = "'29.65"
$r.Style = 'Normal'
$r = $ws.Range('D50')
$r.Value = "'51.67"
$r.Style = 'Normal'
$ws.Range('E50').Value = '  -1.51%  '
$ws.Range('E51').Value = '  +0.11%  '
